$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-02 Monday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-09-03 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("260×7=", $true, $true, $false, $false, $false, $true, 1, $false, "211×2=", 2) | Out-Null
$d.Content.Find.Execute("276×7=", $true, $true, $false, $false, $false, $true, 1, $false, "891×3=", 2) | Out-Null
$d.Content.Find.Execute("444×2=", $true, $true, $false, $false, $false, $true, 1, $false, "437×4=", 2) | Out-Null
$d.Content.Find.Execute("179×7=", $true, $true, $false, $false, $false, $true, 1, $false, "196×9=", 2) | Out-Null
$d.Content.Find.Execute("127×4=", $true, $true, $false, $false, $false, $true, 1, $false, "377×3=", 2) | Out-Null
$d.Content.Find.Execute("569×4=", $true, $true, $false, $false, $false, $true, 1, $false, "716×8=", 2) | Out-Null
$d.Content.Find.Execute("194×8=", $true, $true, $false, $false, $false, $true, 1, $false, "626×4=", 2) | Out-Null
$d.Content.Find.Execute("516×8=", $true, $true, $false, $false, $false, $true, 1, $false, "407×5=", 2) | Out-Null
$d.Content.Find.Execute("502×8=", $true, $true, $false, $false, $false, $true, 1, $false, "104×5=", 2) | Out-Null
$d.Content.Find.Execute("287×5=", $true, $true, $false, $false, $false, $true, 1, $false, "929×8=", 2) | Out-Null
$d.Content.Find.Execute("166×7=", $true, $true, $false, $false, $false, $true, 1, $false, "951×5=", 2) | Out-Null
$d.Content.Find.Execute("569×5=", $true, $true, $false, $false, $false, $true, 1, $false, "449×7=", 2) | Out-Null
$d.Content.Find.Execute("360×3=", $true, $true, $false, $false, $false, $true, 1, $false, "779×3=", 2) | Out-Null
$d.Content.Find.Execute("583×4=", $true, $true, $false, $false, $false, $true, 1, $false, "184×2=", 2) | Out-Null
$d.Content.Find.Execute("857×4=", $true, $true, $false, $false, $false, $true, 1, $false, "497×2=", 2) | Out-Null
$d.Content.Find.Execute("437×8=", $true, $true, $false, $false, $false, $true, 1, $false, "454×4=", 2) | Out-Null
$d.Content.Find.Execute("171×3=", $true, $true, $false, $false, $false, $true, 1, $false, "501×4=", 2) | Out-Null
$d.Content.Find.Execute("455×7=", $true, $true, $false, $false, $false, $true, 1, $false, "345×5=", 2) | Out-Null
$d.Content.Find.Execute("709×8=", $true, $true, $false, $false, $false, $true, 1, $false, "739×8=", 2) | Out-Null
$d.Content.Find.Execute("366×4=", $true, $true, $false, $false, $false, $true, 1, $false, "308×5=", 2) | Out-Null
$d.Content.Find.Execute("133×6=", $true, $true, $false, $false, $false, $true, 1, $false, "669×3=", 2) | Out-Null
$d.Content.Find.Execute("315×3=", $true, $true, $false, $false, $false, $true, 1, $false, "569×7=", 2) | Out-Null
$d.Content.Find.Execute("892×6=", $true, $true, $false, $false, $false, $true, 1, $false, "272×7=", 2) | Out-Null
$d.Content.Find.Execute("220×9=", $true, $true, $false, $false, $false, $true, 1, $false, "466×2=", 2) | Out-Null
$d.Content.Find.Execute("183×4=", $true, $true, $false, $false, $false, $true, 1, $false, "189×3=", 2) | Out-Null
